# Daily update at 8 AM UTC
# Promote the previous "last row" (row 40) to the regular date-time number
# format, then append the new day's data as the new last row (row 41) using
# the date-only number format that always marks the final row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40 is no longer the last row: give it the "interior" number format.
$ws.Range("A40").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New data for 2025-05-04 (serial 45781).
$ws.Range("A41").Value = 45781
$ws.Range("B41").Value = 170
$ws.Range("C41").Value = 175
$ws.Range("D41").Value = 169

# New row 41 becomes the last row: give it the date-only number format.
$ws.Range("A41").NumberFormat = "YYYY-MM-DD"
